function Set-TextCell {
    param($ws, $row, $col, $val)
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextCell $ws 2 4 "68.138.30"
Set-TextCell $ws 2 5 "  +0.12%  "
Set-TextCell $ws 3 4 "3.273.14"
Set-TextCell $ws 3 5 "  +0.54%  "
Set-TextCell $ws 4 5 "  -0.01%  "
Set-TextCell $ws 5 5 "  -0.09%  "
Set-TextCell $ws 6 4 "184.76"
Set-TextCell $ws 6 5 "  +1.87%  "
Set-TextCell $ws 7 5 "  +0.07%  "
Set-TextCell $ws 8 5 "  +1.33%  "
Set-TextCell $ws 9 5 "  -2.99%  "
Set-TextCell $ws 10 5 "  -0.55%  "
Set-TextCell $ws 11 5 "  -2.57%  "
Set-TextCell $ws 12 4 "3.838.26"
Set-TextCell $ws 12 5 "  +0.58%  "
Set-TextCell $ws 13 5 "  +0.98%  "
Set-TextCell $ws 14 5 "  -2.32%  "
Set-TextCell $ws 15 4 "68.142.31"
Set-TextCell $ws 15 5 "  +0.04%  "
Set-TextCell $ws 17 4 "3.267.08"
Set-TextCell $ws 17 5 "  +0.88%  "
Set-TextCell $ws 18 5 "  -1.40%  "
Set-TextCell $ws 19 5 "  -1.09%  "
Set-TextCell $ws 20 4 "418.30"
Set-TextCell $ws 20 5 "  +6.42%  "
Set-TextCell $ws 21 5 "  -1.41%  "
Set-TextCell $ws 22 5 "  -0.05%  "
Set-TextCell $ws 23 4 "71.52"
Set-TextCell $ws 23 5 "  +0.16%  "
Set-TextCell $ws 24 5 "  -1.50%  "
Set-TextCell $ws 25 4 "0.0000118"
Set-TextCell $ws 25 5 "  -1.34%  "
Set-TextCell $ws 26 5 "  -0.95%  "
Set-TextCell $ws 27 4 "9.43"
Set-TextCell $ws 27 5 "  -1.85%  "
Set-TextCell $ws 28 4 "0.998"
Set-TextCell $ws 28 5 "  -0.15%  "
Set-TextCell $ws 29 5 "  -1.88%  "
Set-TextCell $ws 30 4 "22.76"
Set-TextCell $ws 30 5 "  -1.17%  "
Set-TextCell $ws 31 5 "  -3.59%  "
Set-TextCell $ws 32 5 "  -3.07%  "
Set-TextCell $ws 33 5 "  +0.02%  "
Set-TextCell $ws 34 5 "  -2.05%  "
Set-TextCell $ws 35 4 "164.18"
Set-TextCell $ws 35 5 "  -0.40%  "
Set-TextCell $ws 36 5 "  -2.60%  "
Set-TextCell $ws 37 5 "  -1.22%  "
Set-TextCell $ws 38 5 "  +2.98%  "
Set-TextCell $ws 39 5 "  -2.85%  "
Set-TextCell $ws 40 4 "4.48"
Set-TextCell $ws 40 5 "  -2.78%  "
Set-TextCell $ws 41 4 "6.31"
Set-TextCell $ws 41 5 "  -3.83%  "
Set-TextCell $ws 42 4 "2.667.38"
Set-TextCell $ws 42 5 "  +2.75%  "
Set-TextCell $ws 43 5 "  -1.13%  "
Set-TextCell $ws 44 2 "Hedera"
Set-TextCell $ws 44 3 "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell $ws 44 4 "0.0679"
Set-TextCell $ws 44 5 "  -1.56%  "
Set-TextCell $ws 45 2 "dogwifhat"
Set-TextCell $ws 45 3 "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextCell $ws 45 4 "2.44"
Set-TextCell $ws 45 5 "  -1.35%  "
Set-TextCell $ws 46 4 "337.59"
Set-TextCell $ws 46 5 "  -1.48%  "
Set-TextCell $ws 47 4 "24.53"
Set-TextCell $ws 47 5 "  -0.77%  "
Set-TextCell $ws 48 5 "  -2.57%  "
Set-TextCell $ws 49 4 "6.29"
Set-TextCell $ws 49 5 "  -0.09%  "
Set-TextCell $ws 50 4 "0.978"
Set-TextCell $ws 50 5 "  -0.19%  "
Set-TextCell $ws 51 5 "  -1.43%  "
